# Update the "Login Details" sheet: replace the password values in B5 and
# B7 with the names "Cele" and "Nkosi" respectively.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login Details")

$ws.Cells.Item(5, 2).Value = "Cele"
$ws.Cells.Item(7, 2).Value = "Nkosi"

# Reflect the author's on-screen state when the file was saved: zoomed to
# 180% with the view scrolled down and E612 selected.
$ws.Activate()
$ws.Range("E612").Select()
$excel.ActiveWindow.Zoom = 180
$excel.ActiveWindow.ScrollRow = 607
$excel.ActiveWindow.ScrollColumn = 1
